$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 283, shifting existing rows 283-309 down to 286-312
$ws.Range("A283:A285").EntireRow.Insert()

# Row 283
$ws.Range("A283").Value = 7
$ws.Range("B283").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C283").Value = "Ñuble"
$ws.Range("D283").Value = 45106
$ws.Range("E283").Value = 16
$ws.Range("F283").Value = "Fruta"
$ws.Range("G283").Value = 100101
$ws.Range("H283").Value = "Berries"
$ws.Range("I283").Value = 100101007
$ws.Range("J283").Value = "Kiwi"
$ws.Range("K283").Value = "Hayward"
$ws.Range("L283").Value = "Especial"
$ws.Range("M283").Value = 60
$ws.Range("N283").Value = 12000
$ws.Range("O283").Value = 12000
$ws.Range("P283").Value = 12000
$ws.Range("Q283").Value = "$/bandeja 18 kilos"
$ws.Range("R283").Value = "Región de O'Higgins"
$ws.Range("S283").Value = 667
$ws.Range("T283").Value = 18

# Row 284
$ws.Range("A284").Value = 7
$ws.Range("B284").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C284").Value = "Ñuble"
$ws.Range("D284").Value = 45106
$ws.Range("E284").Value = 16
$ws.Range("F284").Value = "Fruta"
$ws.Range("G284").Value = 100101
$ws.Range("H284").Value = "Berries"
$ws.Range("I284").Value = 100101007
$ws.Range("J284").Value = "Kiwi"
$ws.Range("K284").Value = "Hayward"
$ws.Range("L284").Value = "Primera"
$ws.Range("M284").Value = 40
$ws.Range("N284").Value = 10000
$ws.Range("O284").Value = 10000
$ws.Range("P284").Value = 10000
$ws.Range("Q284").Value = "$/bandeja 18 kilos"
$ws.Range("R284").Value = "Región de O'Higgins"
$ws.Range("S284").Value = 556
$ws.Range("T284").Value = 18

# Row 285
$ws.Range("A285").Value = 7
$ws.Range("B285").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C285").Value = "Ñuble"
$ws.Range("D285").Value = 45106
$ws.Range("E285").Value = 16
$ws.Range("F285").Value = "Fruta"
$ws.Range("G285").Value = 100101
$ws.Range("H285").Value = "Berries"
$ws.Range("I285").Value = 100101007
$ws.Range("J285").Value = "Kiwi"
$ws.Range("K285").Value = "Hayward"
$ws.Range("L285").Value = "Segunda"
$ws.Range("M285").Value = 40
$ws.Range("N285").Value = 8000
$ws.Range("O285").Value = 8000
$ws.Range("P285").Value = 8000
$ws.Range("Q285").Value = "$/bandeja 18 kilos"
$ws.Range("R285").Value = "Región de O'Higgins"
$ws.Range("S285").Value = 444
$ws.Range("T285").Value = 18

# Row 286
$ws.Range("A286").Value = 7
$ws.Range("B286").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C286").Value = "Ñuble"
$ws.Range("D286").Value = 44307
$ws.Range("E286").Value = 16
$ws.Range("F286").Value = "Fruta"
$ws.Range("G286").Value = 100101
$ws.Range("H286").Value = "Berries"
$ws.Range("I286").Value = 100101007
$ws.Range("J286").Value = "Kiwi"
$ws.Range("K286").Value = "Hayward"
$ws.Range("L286").Value = "Primera"
$ws.Range("M286").Value = 80
$ws.Range("N286").Value = 10000
$ws.Range("O286").Value = 11000
$ws.Range("P286").Value = 10500
$ws.Range("Q286").Value = "$/bandeja 10 kilos"
$ws.Range("R286").Value = "Provincia de Curicó"
$ws.Range("S286").Value = 1050
$ws.Range("T286").Value = 10

# Row 287
$ws.Range("A287").Value = 7
$ws.Range("B287").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C287").Value = "Ñuble"
$ws.Range("D287").Value = 44307
$ws.Range("E287").Value = 16
$ws.Range("F287").Value = "Fruta"
$ws.Range("G287").Value = 100101
$ws.Range("H287").Value = "Berries"
$ws.Range("I287").Value = 100101007
$ws.Range("J287").Value = "Kiwi"
$ws.Range("K287").Value = "Hayward"
$ws.Range("L287").Value = "Segunda"
$ws.Range("M287").Value = 120
$ws.Range("N287").Value = 8500
$ws.Range("O287").Value = 9000
$ws.Range("P287").Value = 8750
$ws.Range("Q287").Value = "$/bandeja 10 kilos"
$ws.Range("R287").Value = "Provincia de Curicó"
$ws.Range("S287").Value = 875
$ws.Range("T287").Value = 10

# Row 288
$ws.Range("A288").Value = 7
$ws.Range("B288").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C288").Value = "Ñuble"
$ws.Range("D288").Value = 44383
$ws.Range("E288").Value = 16
$ws.Range("F288").Value = "Fruta"
$ws.Range("G288").Value = 100101
$ws.Range("H288").Value = "Berries"
$ws.Range("I288").Value = 100101007
$ws.Range("J288").Value = "Kiwi"
$ws.Range("K288").Value = "Hayward"
$ws.Range("L288").Value = "Primera"
$ws.Range("M288").Value = 120
$ws.Range("N288").Value = 11000
$ws.Range("O288").Value = 12000
$ws.Range("P288").Value = 11500
$ws.Range("Q288").Value = "$/bandeja 18 kilos"
$ws.Range("R288").Value = "Provincia de Curicó"
$ws.Range("S288").Value = 639
$ws.Range("T288").Value = 18

# Row 289
$ws.Range("A289").Value = 7
$ws.Range("B289").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C289").Value = "Ñuble"
$ws.Range("D289").Value = 44383
$ws.Range("E289").Value = 16
$ws.Range("F289").Value = "Fruta"
$ws.Range("G289").Value = 100101
$ws.Range("H289").Value = "Berries"
$ws.Range("I289").Value = 100101007
$ws.Range("J289").Value = "Kiwi"
$ws.Range("K289").Value = "Hayward"
$ws.Range("L289").Value = "Segunda"
$ws.Range("M289").Value = 120
$ws.Range("N289").Value = 9000
$ws.Range("O289").Value = 9500
$ws.Range("P289").Value = 9250
$ws.Range("Q289").Value = "$/bandeja 18 kilos"
$ws.Range("R289").Value = "Provincia de Curicó"
$ws.Range("S289").Value = 514
$ws.Range("T289").Value = 18

# Row 290
$ws.Range("A290").Value = 7
$ws.Range("B290").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C290").Value = "Ñuble"
$ws.Range("D290").Value = 45068
$ws.Range("E290").Value = 16
$ws.Range("F290").Value = "Fruta"
$ws.Range("G290").Value = 100101
$ws.Range("H290").Value = "Berries"
$ws.Range("I290").Value = 100101007
$ws.Range("J290").Value = "Kiwi"
$ws.Range("K290").Value = "Hayward"
$ws.Range("L290").Value = "Especial"
$ws.Range("M290").Value = 50
$ws.Range("N290").Value = 12000
$ws.Range("O290").Value = 12000
$ws.Range("P290").Value = 12000
$ws.Range("Q290").Value = "$/bandeja 18 kilos"
$ws.Range("R290").Value = "Región de O'Higgins"
$ws.Range("S290").Value = 667
$ws.Range("T290").Value = 18

# Row 291
$ws.Range("A291").Value = 7
$ws.Range("B291").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C291").Value = "Ñuble"
$ws.Range("D291").Value = 45068
$ws.Range("E291").Value = 16
$ws.Range("F291").Value = "Fruta"
$ws.Range("G291").Value = 100101
$ws.Range("H291").Value = "Berries"
$ws.Range("I291").Value = 100101007
$ws.Range("J291").Value = "Kiwi"
$ws.Range("K291").Value = "Hayward"
$ws.Range("L291").Value = "Primera"
$ws.Range("M291").Value = 80
$ws.Range("N291").Value = 10000
$ws.Range("O291").Value = 10000
$ws.Range("P291").Value = 10000
$ws.Range("Q291").Value = "$/bandeja 18 kilos"
$ws.Range("R291").Value = "Región de O'Higgins"
$ws.Range("S291").Value = 556
$ws.Range("T291").Value = 18

# Row 292
$ws.Range("A292").Value = 7
$ws.Range("B292").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C292").Value = "Ñuble"
$ws.Range("D292").Value = 44756
$ws.Range("E292").Value = 16
$ws.Range("F292").Value = "Fruta"
$ws.Range("G292").Value = 100101
$ws.Range("H292").Value = "Berries"
$ws.Range("I292").Value = 100101007
$ws.Range("J292").Value = "Kiwi"
$ws.Range("K292").Value = "Hayward"
$ws.Range("L292").Value = "Primera"
$ws.Range("M292").Value = 120
$ws.Range("N292").Value = 6500
$ws.Range("O292").Value = 7000
$ws.Range("P292").Value = 6750
$ws.Range("Q292").Value = "$/bandeja 18 kilos"
$ws.Range("R292").Value = "Provincia de Curicó"
$ws.Range("S292").Value = 375
$ws.Range("T292").Value = 18

# Row 293
$ws.Range("A293").Value = 7
$ws.Range("B293").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C293").Value = "Ñuble"
$ws.Range("D293").Value = 44756
$ws.Range("E293").Value = 16
$ws.Range("F293").Value = "Fruta"
$ws.Range("G293").Value = 100101
$ws.Range("H293").Value = "Berries"
$ws.Range("I293").Value = 100101007
$ws.Range("J293").Value = "Kiwi"
$ws.Range("K293").Value = "Hayward"
$ws.Range("L293").Value = "Segunda"
$ws.Range("M293").Value = 60
$ws.Range("N293").Value = 6000
$ws.Range("O293").Value = 6000
$ws.Range("P293").Value = 6000
$ws.Range("Q293").Value = "$/bandeja 18 kilos"
$ws.Range("R293").Value = "Provincia de Curicó"
$ws.Range("S293").Value = 333
$ws.Range("T293").Value = 18

# Row 294
$ws.Range("A294").Value = 7
$ws.Range("B294").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C294").Value = "Ñuble"
$ws.Range("D294").Value = 44810
$ws.Range("E294").Value = 16
$ws.Range("F294").Value = "Fruta"
$ws.Range("G294").Value = 100101
$ws.Range("H294").Value = "Berries"
$ws.Range("I294").Value = 100101007
$ws.Range("J294").Value = "Kiwi"
$ws.Range("K294").Value = "Hayward"
$ws.Range("L294").Value = "Primera"
$ws.Range("M294").Value = 120
$ws.Range("N294").Value = 8500
$ws.Range("O294").Value = 9000
$ws.Range("P294").Value = 8750
$ws.Range("Q294").Value = "$/bandeja 18 kilos"
$ws.Range("R294").Value = "Provincia de Curicó"
$ws.Range("S294").Value = 486
$ws.Range("T294").Value = 18

# Row 295
$ws.Range("A295").Value = 7
$ws.Range("B295").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C295").Value = "Ñuble"
$ws.Range("D295").Value = 44810
$ws.Range("E295").Value = 16
$ws.Range("F295").Value = "Fruta"
$ws.Range("G295").Value = 100101
$ws.Range("H295").Value = "Berries"
$ws.Range("I295").Value = 100101007
$ws.Range("J295").Value = "Kiwi"
$ws.Range("K295").Value = "Hayward"
$ws.Range("L295").Value = "Segunda"
$ws.Range("M295").Value = 60
$ws.Range("N295").Value = 7000
$ws.Range("O295").Value = 7000
$ws.Range("P295").Value = 7000
$ws.Range("Q295").Value = "$/bandeja 18 kilos"
$ws.Range("R295").Value = "Provincia de Curicó"
$ws.Range("S295").Value = 389
$ws.Range("T295").Value = 18

# Row 296
$ws.Range("A296").Value = 7
$ws.Range("B296").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C296").Value = "Ñuble"
$ws.Range("D296").Value = 44442
$ws.Range("E296").Value = 16
$ws.Range("F296").Value = "Fruta"
$ws.Range("G296").Value = 100101
$ws.Range("H296").Value = "Berries"
$ws.Range("I296").Value = 100101007
$ws.Range("J296").Value = "Kiwi"
$ws.Range("K296").Value = "Hayward"
$ws.Range("L296").Value = "Primera"
$ws.Range("M296").Value = 120
$ws.Range("N296").Value = 12000
$ws.Range("O296").Value = 12500
$ws.Range("P296").Value = 12250
$ws.Range("Q296").Value = "$/bandeja 18 kilos"
$ws.Range("R296").Value = "Provincia de Curicó"
$ws.Range("S296").Value = 681
$ws.Range("T296").Value = 18

# Row 297
$ws.Range("A297").Value = 7
$ws.Range("B297").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C297").Value = "Ñuble"
$ws.Range("D297").Value = 44442
$ws.Range("E297").Value = 16
$ws.Range("F297").Value = "Fruta"
$ws.Range("G297").Value = 100101
$ws.Range("H297").Value = "Berries"
$ws.Range("I297").Value = 100101007
$ws.Range("J297").Value = "Kiwi"
$ws.Range("K297").Value = "Hayward"
$ws.Range("L297").Value = "Segunda"
$ws.Range("M297").Value = 100
$ws.Range("N297").Value = 11000
$ws.Range("O297").Value = 11500
$ws.Range("P297").Value = 11250
$ws.Range("Q297").Value = "$/bandeja 18 kilos"
$ws.Range("R297").Value = "Provincia de Curicó"
$ws.Range("S297").Value = 625
$ws.Range("T297").Value = 18

# Row 298
$ws.Range("A298").Value = 7
$ws.Range("B298").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C298").Value = "Ñuble"
$ws.Range("D298").Value = 44704
$ws.Range("E298").Value = 16
$ws.Range("F298").Value = "Fruta"
$ws.Range("G298").Value = 100101
$ws.Range("H298").Value = "Berries"
$ws.Range("I298").Value = 100101007
$ws.Range("J298").Value = "Kiwi"
$ws.Range("K298").Value = "Hayward"
$ws.Range("L298").Value = "Primera"
$ws.Range("M298").Value = 120
$ws.Range("N298").Value = 11000
$ws.Range("O298").Value = 12000
$ws.Range("P298").Value = 11500
$ws.Range("Q298").Value = "$/bandeja 18 kilos"
$ws.Range("R298").Value = "Provincia de Curicó"
$ws.Range("S298").Value = 639
$ws.Range("T298").Value = 18

# Row 299
$ws.Range("A299").Value = 7
$ws.Range("B299").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C299").Value = "Ñuble"
$ws.Range("D299").Value = 44704
$ws.Range("E299").Value = 16
$ws.Range("F299").Value = "Fruta"
$ws.Range("G299").Value = 100101
$ws.Range("H299").Value = "Berries"
$ws.Range("I299").Value = 100101007
$ws.Range("J299").Value = "Kiwi"
$ws.Range("K299").Value = "Hayward"
$ws.Range("L299").Value = "Segunda"
$ws.Range("M299").Value = 120
$ws.Range("N299").Value = 8500
$ws.Range("O299").Value = 9000
$ws.Range("P299").Value = 8750
$ws.Range("Q299").Value = "$/bandeja 18 kilos"
$ws.Range("R299").Value = "Provincia de Curicó"
$ws.Range("S299").Value = 486
$ws.Range("T299").Value = 18

# Row 300
$ws.Range("A300").Value = 7
$ws.Range("B300").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C300").Value = "Ñuble"
$ws.Range("D300").Value = 44769
$ws.Range("E300").Value = 16
$ws.Range("F300").Value = "Fruta"
$ws.Range("G300").Value = 100101
$ws.Range("H300").Value = "Berries"
$ws.Range("I300").Value = 100101007
$ws.Range("J300").Value = "Kiwi"
$ws.Range("K300").Value = "Hayward"
$ws.Range("L300").Value = "Primera"
$ws.Range("M300").Value = 120
$ws.Range("N300").Value = 6500
$ws.Range("O300").Value = 7000
$ws.Range("P300").Value = 6750
$ws.Range("Q300").Value = "$/bandeja 18 kilos"
$ws.Range("R300").Value = "Provincia de Curicó"
$ws.Range("S300").Value = 375
$ws.Range("T300").Value = 18

# Row 301
$ws.Range("A301").Value = 7
$ws.Range("B301").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C301").Value = "Ñuble"
$ws.Range("D301").Value = 44769
$ws.Range("E301").Value = 16
$ws.Range("F301").Value = "Fruta"
$ws.Range("G301").Value = 100101
$ws.Range("H301").Value = "Berries"
$ws.Range("I301").Value = 100101007
$ws.Range("J301").Value = "Kiwi"
$ws.Range("K301").Value = "Hayward"
$ws.Range("L301").Value = "Segunda"
$ws.Range("M301").Value = 120
$ws.Range("N301").Value = 5500
$ws.Range("O301").Value = 6000
$ws.Range("P301").Value = 5750
$ws.Range("Q301").Value = "$/bandeja 18 kilos"
$ws.Range("R301").Value = "Provincia de Curicó"
$ws.Range("S301").Value = 319
$ws.Range("T301").Value = 18

# Row 302
$ws.Range("A302").Value = 7
$ws.Range("B302").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C302").Value = "Ñuble"
$ws.Range("D302").Value = 45104
$ws.Range("E302").Value = 16
$ws.Range("F302").Value = "Fruta"
$ws.Range("G302").Value = 100101
$ws.Range("H302").Value = "Berries"
$ws.Range("I302").Value = 100101007
$ws.Range("J302").Value = "Kiwi"
$ws.Range("K302").Value = "Hayward"
$ws.Range("L302").Value = "Especial"
$ws.Range("M302").Value = 60
$ws.Range("N302").Value = 12000
$ws.Range("O302").Value = 12000
$ws.Range("P302").Value = 12000
$ws.Range("Q302").Value = "$/bandeja 18 kilos"
$ws.Range("R302").Value = "Región de O'Higgins"
$ws.Range("S302").Value = 667
$ws.Range("T302").Value = 18

# Row 303
$ws.Range("A303").Value = 7
$ws.Range("B303").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C303").Value = "Ñuble"
$ws.Range("D303").Value = 45104
$ws.Range("E303").Value = 16
$ws.Range("F303").Value = "Fruta"
$ws.Range("G303").Value = 100101
$ws.Range("H303").Value = "Berries"
$ws.Range("I303").Value = 100101007
$ws.Range("J303").Value = "Kiwi"
$ws.Range("K303").Value = "Hayward"
$ws.Range("L303").Value = "Primera"
$ws.Range("M303").Value = 40
$ws.Range("N303").Value = 10000
$ws.Range("O303").Value = 10000
$ws.Range("P303").Value = 10000
$ws.Range("Q303").Value = "$/bandeja 18 kilos"
$ws.Range("R303").Value = "Región de O'Higgins"
$ws.Range("S303").Value = 556
$ws.Range("T303").Value = 18

# Row 304
$ws.Range("A304").Value = 7
$ws.Range("B304").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C304").Value = "Ñuble"
$ws.Range("D304").Value = 45104
$ws.Range("E304").Value = 16
$ws.Range("F304").Value = "Fruta"
$ws.Range("G304").Value = 100101
$ws.Range("H304").Value = "Berries"
$ws.Range("I304").Value = 100101007
$ws.Range("J304").Value = "Kiwi"
$ws.Range("K304").Value = "Hayward"
$ws.Range("L304").Value = "Segunda"
$ws.Range("M304").Value = 30
$ws.Range("N304").Value = 8000
$ws.Range("O304").Value = 8000
$ws.Range("P304").Value = 8000
$ws.Range("Q304").Value = "$/bandeja 18 kilos"
$ws.Range("R304").Value = "Región de O'Higgins"
$ws.Range("S304").Value = 444
$ws.Range("T304").Value = 18

# Row 305
$ws.Range("A305").Value = 7
$ws.Range("B305").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C305").Value = "Ñuble"
$ws.Range("D305").Value = 45096
$ws.Range("E305").Value = 16
$ws.Range("F305").Value = "Fruta"
$ws.Range("G305").Value = 100101
$ws.Range("H305").Value = "Berries"
$ws.Range("I305").Value = 100101007
$ws.Range("J305").Value = "Kiwi"
$ws.Range("K305").Value = "Hayward"
$ws.Range("L305").Value = "Especial"
$ws.Range("M305").Value = 60
$ws.Range("N305").Value = 12000
$ws.Range("O305").Value = 12000
$ws.Range("P305").Value = 12000
$ws.Range("Q305").Value = "$/bandeja 18 kilos"
$ws.Range("R305").Value = "Región de O'Higgins"
$ws.Range("S305").Value = 667
$ws.Range("T305").Value = 18

# Row 306
$ws.Range("A306").Value = 7
$ws.Range("B306").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C306").Value = "Ñuble"
$ws.Range("D306").Value = 45096
$ws.Range("E306").Value = 16
$ws.Range("F306").Value = "Fruta"
$ws.Range("G306").Value = 100101
$ws.Range("H306").Value = "Berries"
$ws.Range("I306").Value = 100101007
$ws.Range("J306").Value = "Kiwi"
$ws.Range("K306").Value = "Hayward"
$ws.Range("L306").Value = "Primera"
$ws.Range("M306").Value = 60
$ws.Range("N306").Value = 10000
$ws.Range("O306").Value = 10000
$ws.Range("P306").Value = 10000
$ws.Range("Q306").Value = "$/bandeja 18 kilos"
$ws.Range("R306").Value = "Región de O'Higgins"
$ws.Range("S306").Value = 556
$ws.Range("T306").Value = 18

# Row 307
$ws.Range("A307").Value = 7
$ws.Range("B307").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C307").Value = "Ñuble"
$ws.Range("D307").Value = 44449
$ws.Range("E307").Value = 16
$ws.Range("F307").Value = "Fruta"
$ws.Range("G307").Value = 100101
$ws.Range("H307").Value = "Berries"
$ws.Range("I307").Value = 100101007
$ws.Range("J307").Value = "Kiwi"
$ws.Range("K307").Value = "Hayward"
$ws.Range("L307").Value = "Primera"
$ws.Range("M307").Value = 100
$ws.Range("N307").Value = 12000
$ws.Range("O307").Value = 12500
$ws.Range("P307").Value = 12250
$ws.Range("Q307").Value = "$/bandeja 18 kilos"
$ws.Range("R307").Value = "Provincia de Curicó"
$ws.Range("S307").Value = 681
$ws.Range("T307").Value = 18

# Row 308
$ws.Range("A308").Value = 7
$ws.Range("B308").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C308").Value = "Ñuble"
$ws.Range("D308").Value = 44449
$ws.Range("E308").Value = 16
$ws.Range("F308").Value = "Fruta"
$ws.Range("G308").Value = 100101
$ws.Range("H308").Value = "Berries"
$ws.Range("I308").Value = 100101007
$ws.Range("J308").Value = "Kiwi"
$ws.Range("K308").Value = "Hayward"
$ws.Range("L308").Value = "Segunda"
$ws.Range("M308").Value = 60
$ws.Range("N308").Value = 11000
$ws.Range("O308").Value = 11500
$ws.Range("P308").Value = 11250
$ws.Range("Q308").Value = "$/bandeja 18 kilos"
$ws.Range("R308").Value = "Provincia de Curicó"
$ws.Range("S308").Value = 625
$ws.Range("T308").Value = 18

# Row 309
$ws.Range("A309").Value = 7
$ws.Range("B309").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C309").Value = "Ñuble"
$ws.Range("D309").Value = 45054
$ws.Range("E309").Value = 16
$ws.Range("F309").Value = "Fruta"
$ws.Range("G309").Value = 100101
$ws.Range("H309").Value = "Berries"
$ws.Range("I309").Value = 100101007
$ws.Range("J309").Value = "Kiwi"
$ws.Range("K309").Value = "Hayward"
$ws.Range("L309").Value = "Primera"
$ws.Range("M309").Value = 100
$ws.Range("N309").Value = 12000
$ws.Range("O309").Value = 13000
$ws.Range("P309").Value = 12500
$ws.Range("Q309").Value = "$/bandeja 18 kilos"
$ws.Range("R309").Value = "Región de O'Higgins"
$ws.Range("S309").Value = 694
$ws.Range("T309").Value = 18

# Row 310
$ws.Range("A310").Value = 7
$ws.Range("B310").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C310").Value = "Ñuble"
$ws.Range("D310").Value = 45049
$ws.Range("E310").Value = 16
$ws.Range("F310").Value = "Fruta"
$ws.Range("G310").Value = 100101
$ws.Range("H310").Value = "Berries"
$ws.Range("I310").Value = 100101007
$ws.Range("J310").Value = "Kiwi"
$ws.Range("K310").Value = "Hayward"
$ws.Range("L310").Value = "Primera"
$ws.Range("M310").Value = 80
$ws.Range("N310").Value = 13000
$ws.Range("O310").Value = 13000
$ws.Range("P310").Value = 13000
$ws.Range("Q310").Value = "$/bandeja 18 kilos"
$ws.Range("R310").Value = "Región de O'Higgins"
$ws.Range("S310").Value = 722
$ws.Range("T310").Value = 18

# Row 311
$ws.Range("A311").Value = 7
$ws.Range("B311").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C311").Value = "Ñuble"
$ws.Range("D311").Value = 44455
$ws.Range("E311").Value = 16
$ws.Range("F311").Value = "Fruta"
$ws.Range("G311").Value = 100101
$ws.Range("H311").Value = "Berries"
$ws.Range("I311").Value = 100101007
$ws.Range("J311").Value = "Kiwi"
$ws.Range("K311").Value = "Hayward"
$ws.Range("L311").Value = "Primera"
$ws.Range("M311").Value = 120
$ws.Range("N311").Value = 12000
$ws.Range("O311").Value = 12500
$ws.Range("P311").Value = 12250
$ws.Range("Q311").Value = "$/bandeja 18 kilos"
$ws.Range("R311").Value = "Provincia de Curicó"
$ws.Range("S311").Value = 681
$ws.Range("T311").Value = 18

# Row 312
$ws.Range("A312").Value = 7
$ws.Range("B312").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C312").Value = "Ñuble"
$ws.Range("D312").Value = 44455
$ws.Range("E312").Value = 16
$ws.Range("F312").Value = "Fruta"
$ws.Range("G312").Value = 100101
$ws.Range("H312").Value = "Berries"
$ws.Range("I312").Value = 100101007
$ws.Range("J312").Value = "Kiwi"
$ws.Range("K312").Value = "Hayward"
$ws.Range("L312").Value = "Segunda"
$ws.Range("M312").Value = 60
$ws.Range("N312").Value = 11000
$ws.Range("O312").Value = 11500
$ws.Range("P312").Value = 11250
$ws.Range("Q312").Value = "$/bandeja 18 kilos"
$ws.Range("R312").Value = "Provincia de Curicó"
$ws.Range("S312").Value = 625
$ws.Range("T312").Value = 18
